{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" problems in the worksheet table\n// with their updated values, one search/replace per distinct original\n// expression (each original string occurs exactly once in the document).\nconst replacements = [\n  [\"188\u00f76=\", \"976\u00f77=\"],\n  [\"712\u00f74=\", \"510\u00f74=\"],\n  [\"123\u00f78=\", \"582\u00f75=\"],\n  [\"426\u00f73=\", \"850\u00f76=\"],\n  [\"927\u00f79=\", \"596\u00f78=\"],\n  [\"833\u00f75=\", \"542\u00f74=\"],\n  [\"527\u00f73=\", \"651\u00f75=\"],\n  [\"279\u00f76=\", \"397\u00f72=\"],\n  [\"985\u00f75=\", \"851\u00f79=\"],\n  [\"470\u00f73=\", \"483\u00f75=\"],\n  [\"491\u00f74=\", \"408\u00f74=\"],\n  [\"788\u00f78=\", \"506\u00f73=\"],\n  [\"912\u00f77=\", \"600\u00f72=\"],\n  [\"194\u00f72=\", \"110\u00f75=\"],\n  [\"772\u00f74=\", \"913\u00f76=\"],\n  [\"558\u00f79=\", \"898\u00f79=\"],\n  [\"545\u00f75=\", \"595\u00f73=\"],\n  [\"609\u00f78=\", \"380\u00f73=\"],\n  [\"285\u00f73=\", \"137\u00f72=\"],\n  [\"413\u00f73=\", \"612\u00f76=\"],\n  [\"150\u00f78=\", \"140\u00f73=\"],\n  [\"116\u00f77=\", \"424\u00f78=\"],\n  [\"781\u00f73=\", \"178\u00f75=\"],\n  [\"348\u00f77=\", \"176\u00f79=\"],\n  [\"830\u00f78=\", \"262\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" problems in the worksheet table\n# with their updated values, one Find/Replace per distinct original\n# expression (each original string occurs exactly once in the document).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"188\u00f76=\", \"976\u00f77=\"),\n    @(\"712\u00f74=\", \"510\u00f74=\"),\n    @(\"123\u00f78=\", \"582\u00f75=\"),\n    @(\"426\u00f73=\", \"850\u00f76=\"),\n    @(\"927\u00f79=\", \"596\u00f78=\"),\n    @(\"833\u00f75=\", \"542\u00f74=\"),\n    @(\"527\u00f73=\", \"651\u00f75=\"),\n    @(\"279\u00f76=\", \"397\u00f72=\"),\n    @(\"985\u00f75=\", \"851\u00f79=\"),\n    @(\"470\u00f73=\", \"483\u00f75=\"),\n    @(\"491\u00f74=\", \"408\u00f74=\"),\n    @(\"788\u00f78=\", \"506\u00f73=\"),\n    @(\"912\u00f77=\", \"600\u00f72=\"),\n    @(\"194\u00f72=\", \"110\u00f75=\"),\n    @(\"772\u00f74=\", \"913\u00f76=\"),\n    @(\"558\u00f79=\", \"898\u00f79=\"),\n    @(\"545\u00f75=\", \"595\u00f73=\"),\n    @(\"609\u00f78=\", \"380\u00f73=\"),\n    @(\"285\u00f73=\", \"137\u00f72=\"),\n    @(\"413\u00f73=\", \"612\u00f76=\"),\n    @(\"150\u00f78=\", \"140\u00f73=\"),\n    @(\"116\u00f77=\", \"424\u00f78=\"),\n    @(\"781\u00f73=\", \"178\u00f75=\"),\n    @(\"348\u00f77=\", \"176\u00f79=\"),\n    @(\"830\u00f78=\", \"262\u00f75=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.Text = $to\n    $find.Execute($from, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $to, $wdReplaceAll) | Out-Null\n}\n"}
